$d = $word.ActiveDocument

# Locate the paragraph that begins the trailing site-chrome block:
#   "Ver no Jupiter Salvar em pdf Salvar em docx"
# Per the diff, this paragraph, the copyright/footer paragraph right after it
# ("(c) 2020 . Contact: luizeleno@usp.br. ...") and the blank paragraph that
# follows both of them are all removed, while the "LOB1004: Calculo II
# (Requisito fraco)" paragraph and its own trailing blank paragraph (as well
# as the page-break paragraph further down) are left untouched.
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Ver no Jupiter*") {
        $target = $p
        break
    }
}

if ($target -ne $null) {
    $start = $target.Range.Start

    # The block to delete spans this paragraph plus the two that follow it.
    $p2 = $target.Next()
    $p3 = $p2.Next()
    $end = $p3.Range.End

    $d.Range($start, $end).Delete()
}
